$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking values stored as text (column D) need NumberFormat "@"
# forced before assignment, then ClearFormats() to strip the temporary
# text-format style so the cell keeps its original (default) style index
# while the underlying value remains a text string, matching the source
# workbook where these Price cells are inline strings, not numbers.
function Set-TextValue($ws, $addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws "D2" "247.53"
Set-TextValue $ws "D3" "22.78"
Set-TextValue $ws "D4" "5.290"
Set-TextValue $ws "D6" "3.426"
Set-TextValue $ws "D8" "0.8100"
Set-TextValue $ws "D9" "0.8664"
Set-TextValue $ws "D10" "0.1433"
Set-TextValue $ws "D11" "0.07389"
Set-TextValue $ws "D12" "0.03049"
Set-TextValue $ws "D14" "0.09381"
Set-TextValue $ws "D15" "3.889"
Set-TextValue $ws "D16" "0.001580"
Set-TextValue $ws "D17" "0.04822"
Set-TextValue $ws "D18" "0.0005849"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue $ws "D20" "0.005170"
Set-TextValue $ws "D21" "0.0009963"
Set-TextValue $ws "D23" "3.737"
Set-TextValue $ws "D24" "2.200"
Set-TextValue $ws "D25" "0.3276"
Set-TextValue $ws "D26" "0.1296"
Set-TextValue $ws "D40" "0.03947"
Set-TextValue $ws "D41" "0.006750"
Set-TextValue $ws "D43" "0.003200"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
Set-TextValue $ws "D44" "0.007993"
$ws.Range("E44").Value = "43LocalTradersLCT"
Set-TextValue $ws "D45" "0.00005622"
Set-TextValue $ws "D47" "0.3599"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
Set-TextValue $ws "D48" "0.1807"
